$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the three new rows (shift subsequent rows down)
$ws.Rows.Item(15).Insert()   # new row for "Glory"
$ws.Rows.Item(20).Insert()   # new row for "Entomb"
$ws.Rows.Item(66).Insert()   # new row for "Helm of Awakening"

# Match the row height used by the surrounding data rows
$ws.Rows.Item(15).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(66).RowHeight = 15

# Update Name / Set / Version / Price for every data row (2-69)
$ws.Cells.Item(2, 1).Value = 'Godless Shrine'
$ws.Cells.Item(2, 2).Value = 'Gatecrash'
$ws.Cells.Item(2, 3).Value = 'Normal'
$ws.Cells.Item(2, 4).Value = 9.7
$ws.Cells.Item(3, 1).Value = 'Morphic Pool'
$ws.Cells.Item(3, 2).Value = 'Commander Legends: Battle For Baldur''s Gate'
$ws.Cells.Item(3, 3).Value = 'Normal'
$ws.Cells.Item(3, 4).Value = 8.72
$ws.Cells.Item(4, 1).Value = 'Sacred Foundry'
$ws.Cells.Item(4, 2).Value = 'Gatecrash'
$ws.Cells.Item(4, 3).Value = 'Normal'
$ws.Cells.Item(4, 4).Value = 13.29
$ws.Cells.Item(5, 1).Value = 'Blood Crypt'
$ws.Cells.Item(5, 2).Value = 'Return to Ravnica'
$ws.Cells.Item(5, 3).Value = 'Normal'
$ws.Cells.Item(5, 4).Value = 16.32
$ws.Cells.Item(6, 1).Value = 'Dragonskull Summit'
$ws.Cells.Item(6, 2).Value = 'Magic 2013'
$ws.Cells.Item(6, 3).Value = 'Normal'
$ws.Cells.Item(6, 4).Value = 2.94
$ws.Cells.Item(7, 1).Value = 'Stomping Ground'
$ws.Cells.Item(7, 2).Value = 'Gatecrash'
$ws.Cells.Item(7, 3).Value = 'Normal'
$ws.Cells.Item(7, 4).Value = 8.72
$ws.Cells.Item(8, 1).Value = 'Jetmir''s Garden'
$ws.Cells.Item(8, 2).Value = 'Streets of New Capenna'
$ws.Cells.Item(8, 3).Value = 'Normal'
$ws.Cells.Item(8, 4).Value = 7.95
$ws.Cells.Item(9, 1).Value = 'Maze''s End'
$ws.Cells.Item(9, 2).Value = 'Dragon''s Maze'
$ws.Cells.Item(9, 3).Value = 'Normal'
$ws.Cells.Item(9, 4).Value = 4.32
$ws.Cells.Item(10, 1).Value = 'Nykthos, Shrine to Nyx'
$ws.Cells.Item(10, 2).Value = 'Theros'
$ws.Cells.Item(10, 3).Value = 'Normal'
$ws.Cells.Item(10, 4).Value = 26.63
$ws.Cells.Item(11, 1).Value = 'Yavimaya Coast'
$ws.Cells.Item(11, 2).Value = 'Dominaria United: Extras'
$ws.Cells.Item(11, 3).Value = 'Normal'
$ws.Cells.Item(11, 4).Value = 1.47
$ws.Cells.Item(12, 1).Value = 'Argoth, Sanctum of Nature: Titania, Gaea Incarnate'
$ws.Cells.Item(12, 2).Value = 'The Brothers'' War Promos'
$ws.Cells.Item(12, 3).Value = 'Normal'
$ws.Cells.Item(12, 4).Value = 4.38
$ws.Cells.Item(13, 1).Value = 'Lair of the Hydra'
$ws.Cells.Item(13, 2).Value = 'Adventures in the Forgotten Realms'
$ws.Cells.Item(13, 3).Value = 'Normal'
$ws.Cells.Item(13, 4).Value = 1.31
$ws.Cells.Item(14, 1).Value = 'Ajani Steadfast'
$ws.Cells.Item(14, 2).Value = 'Magic 2015'
$ws.Cells.Item(14, 3).Value = 'Normal'
$ws.Cells.Item(14, 4).Value = 6.11
$ws.Cells.Item(15, 1).Value = 'Glory'
$ws.Cells.Item(15, 2).Value = 'Dominaria Remastered: Extras'
$ws.Cells.Item(15, 3).Value = 'Foil'
$ws.Cells.Item(15, 4).Value = 1.44
$ws.Cells.Item(16, 1).Value = 'Sage of Hours'
$ws.Cells.Item(16, 2).Value = 'Journey into Nyx'
$ws.Cells.Item(16, 3).Value = 'Normal'
$ws.Cells.Item(16, 4).Value = 3.95
$ws.Cells.Item(17, 1).Value = 'Ancient Brass Dragon'
$ws.Cells.Item(17, 2).Value = 'Commander Legends: Battle For Baldur''s Gate'
$ws.Cells.Item(17, 3).Value = 'Normal'
$ws.Cells.Item(17, 4).Value = 14.1
$ws.Cells.Item(18, 1).Value = 'Dark Petition'
$ws.Cells.Item(18, 2).Value = 'Magic Origins'
$ws.Cells.Item(18, 3).Value = 'Normal'
$ws.Cells.Item(18, 4).Value = 3.89
$ws.Cells.Item(19, 1).Value = 'Diabolic Intent'
$ws.Cells.Item(19, 2).Value = 'The Brothers'' War'
$ws.Cells.Item(19, 3).Value = 'Normal'
$ws.Cells.Item(19, 4).Value = 6.93
$ws.Cells.Item(20, 1).Value = 'Entomb'
$ws.Cells.Item(20, 2).Value = 'Dominaria Remastered: Extras'
$ws.Cells.Item(20, 3).Value = 'V.2'
$ws.Cells.Item(20, 4).Value = 19.83
$ws.Cells.Item(21, 1).Value = 'Gravecrawler'
$ws.Cells.Item(21, 2).Value = 'Dark Ascension'
$ws.Cells.Item(21, 3).Value = 'Normal'
$ws.Cells.Item(21, 4).Value = 3.48
$ws.Cells.Item(22, 1).Value = 'Scourge of the Skyclaves'
$ws.Cells.Item(22, 2).Value = 'Zendikar Rising Promos'
$ws.Cells.Item(22, 3).Value = 'V.1'
$ws.Cells.Item(22, 4).Value = 0.74
$ws.Cells.Item(23, 1).Value = 'Tree of Perdition'
$ws.Cells.Item(23, 2).Value = 'Eldritch Moon'
$ws.Cells.Item(23, 3).Value = 'Normal'
$ws.Cells.Item(23, 4).Value = 5.21
$ws.Cells.Item(24, 1).Value = 'Deadly Dispute'
$ws.Cells.Item(24, 2).Value = 'Commander Legends: Battle For Baldur''s Gate'
$ws.Cells.Item(24, 3).Value = 'Foil'
$ws.Cells.Item(24, 4).Value = 2.36
$ws.Cells.Item(25, 1).Value = 'Brotherhood''s End'
$ws.Cells.Item(25, 2).Value = 'The Brothers'' War'
$ws.Cells.Item(25, 3).Value = 'Normal'
$ws.Cells.Item(25, 4).Value = 7.32
$ws.Cells.Item(26, 1).Value = 'Flame-Wreathed Phoenix'
$ws.Cells.Item(26, 2).Value = 'Born of the Gods'
$ws.Cells.Item(26, 3).Value = 'Foil'
$ws.Cells.Item(26, 4).Value = 1.58
$ws.Cells.Item(27, 1).Value = 'Shivan Devastator'
$ws.Cells.Item(27, 2).Value = 'Dominaria United'
$ws.Cells.Item(27, 3).Value = 'Normal'
$ws.Cells.Item(27, 4).Value = 5.05
$ws.Cells.Item(28, 1).Value = 'Vexing Devil'
$ws.Cells.Item(28, 2).Value = 'Avacyn Restored'
$ws.Cells.Item(28, 3).Value = 'Normal'
$ws.Cells.Item(28, 4).Value = 4.16
$ws.Cells.Item(29, 1).Value = 'Fauna Shaman'
$ws.Cells.Item(29, 2).Value = 'The Brothers'' War'
$ws.Cells.Item(29, 3).Value = 'Normal'
$ws.Cells.Item(29, 4).Value = 1.06
$ws.Cells.Item(30, 1).Value = 'Fyndhorn Elves'
$ws.Cells.Item(30, 2).Value = '30th Anniversary Celebration'
$ws.Cells.Item(30, 3).Value = 'German'
$ws.Cells.Item(30, 4).Value = 1.7
$ws.Cells.Item(31, 1).Value = 'Gyre Sage'
$ws.Cells.Item(31, 2).Value = 'Gatecrash'
$ws.Cells.Item(31, 3).Value = 'Normal'
$ws.Cells.Item(31, 4).Value = 1.4
$ws.Cells.Item(32, 1).Value = 'Hardened Scales'
$ws.Cells.Item(32, 2).Value = 'Khans of Tarkir'
$ws.Cells.Item(32, 3).Value = 'Normal'
$ws.Cells.Item(32, 4).Value = 2.88
$ws.Cells.Item(33, 1).Value = 'Life from the Loam'
$ws.Cells.Item(33, 2).Value = 'Duel Decks: Izzet vs Golgari'
$ws.Cells.Item(33, 3).Value = 'Normal'
$ws.Cells.Item(33, 4).Value = 10.6
$ws.Cells.Item(34, 1).Value = 'Majestic Genesis'
$ws.Cells.Item(34, 2).Value = 'Commander Legends: Battle For Baldur''s Gate Promos'
$ws.Cells.Item(34, 3).Value = 'Foil'
$ws.Cells.Item(34, 4).Value = 2.9
$ws.Cells.Item(35, 1).Value = 'Majestic Genesis'
$ws.Cells.Item(35, 2).Value = 'Commander Legends: Battle For Baldur''s Gate'
$ws.Cells.Item(35, 3).Value = 'Normal'
$ws.Cells.Item(35, 4).Value = 1.65
$ws.Cells.Item(36, 1).Value = 'Nylea, God of the Hunt'
$ws.Cells.Item(36, 2).Value = 'Theros'
$ws.Cells.Item(36, 3).Value = 'Normal'
$ws.Cells.Item(36, 4).Value = 4.93
$ws.Cells.Item(37, 1).Value = 'Parallel Lives'
$ws.Cells.Item(37, 2).Value = 'Innistrad'
$ws.Cells.Item(37, 3).Value = 'Normal'
$ws.Cells.Item(37, 4).Value = 29.65
$ws.Cells.Item(38, 1).Value = 'Silverback Elder'
$ws.Cells.Item(38, 2).Value = 'Dominaria United: Extras'
$ws.Cells.Item(38, 3).Value = 'Normal'
$ws.Cells.Item(38, 4).Value = 9.73
$ws.Cells.Item(39, 1).Value = 'Haywire Mite'
$ws.Cells.Item(39, 2).Value = 'The Brothers'' War'
$ws.Cells.Item(39, 3).Value = 'Normal'
$ws.Cells.Item(39, 4).Value = 1.23
$ws.Cells.Item(40, 1).Value = 'Sphinx''s Revelation'
$ws.Cells.Item(40, 2).Value = 'Return to Ravnica'
$ws.Cells.Item(40, 3).Value = 'Normal'
$ws.Cells.Item(40, 4).Value = 2.05
$ws.Cells.Item(41, 1).Value = 'Ashiok, Nightmare Weaver'
$ws.Cells.Item(41, 2).Value = 'Theros'
$ws.Cells.Item(41, 3).Value = 'Normal'
$ws.Cells.Item(41, 4).Value = 3.47
$ws.Cells.Item(42, 1).Value = 'Mind Grind'
$ws.Cells.Item(42, 2).Value = 'Gatecrash'
$ws.Cells.Item(42, 3).Value = 'Normal'
$ws.Cells.Item(42, 4).Value = 3.79
$ws.Cells.Item(43, 1).Value = 'Satoru Umezawa'
$ws.Cells.Item(43, 2).Value = 'Buy a Box Promos'
$ws.Cells.Item(43, 3).Value = 'Normal'
$ws.Cells.Item(43, 4).Value = 0.42
$ws.Cells.Item(44, 1).Value = 'Legion''s Initiative'
$ws.Cells.Item(44, 2).Value = 'Dragon''s Maze'
$ws.Cells.Item(44, 3).Value = 'Normal'
$ws.Cells.Item(44, 4).Value = 1.37
$ws.Cells.Item(45, 1).Value = 'Expressive Iteration'
$ws.Cells.Item(45, 2).Value = 'Strixhaven: School of Mages'
$ws.Cells.Item(45, 3).Value = 'Normal'
$ws.Cells.Item(45, 4).Value = 3.19
$ws.Cells.Item(46, 1).Value = 'Jhoira, Ageless Innovator'
$ws.Cells.Item(46, 2).Value = 'Dominaria United: Extras'
$ws.Cells.Item(46, 3).Value = 'V.2'
$ws.Cells.Item(46, 4).Value = 1.46
$ws.Cells.Item(47, 1).Value = 'Kolaghan''s Command'
$ws.Cells.Item(47, 2).Value = 'Double Masters 2022'
$ws.Cells.Item(47, 3).Value = 'Normal'
$ws.Cells.Item(47, 4).Value = 2.99
$ws.Cells.Item(48, 1).Value = 'Mogis, God of Slaughter'
$ws.Cells.Item(48, 2).Value = 'Born of the Gods'
$ws.Cells.Item(48, 3).Value = 'Normal'
$ws.Cells.Item(48, 4).Value = 7.53
$ws.Cells.Item(49, 1).Value = 'Ajani, Sleeper Agent'
$ws.Cells.Item(49, 2).Value = 'Dominaria United'
$ws.Cells.Item(49, 3).Value = 'Foil'
$ws.Cells.Item(49, 4).Value = 3.37
$ws.Cells.Item(50, 1).Value = 'Ivy, Gleeful Spellthief'
$ws.Cells.Item(50, 2).Value = 'Dominaria United Promos'
$ws.Cells.Item(50, 3).Value = 'V.1'
$ws.Cells.Item(50, 4).Value = 1.47
$ws.Cells.Item(51, 1).Value = 'Vorel of the Hull Clade'
$ws.Cells.Item(51, 2).Value = 'Dragon''s Maze'
$ws.Cells.Item(51, 3).Value = 'Foil'
$ws.Cells.Item(51, 4).Value = 1.66
$ws.Cells.Item(52, 1).Value = 'Deathrite Shaman'
$ws.Cells.Item(52, 2).Value = 'Return to Ravnica'
$ws.Cells.Item(52, 3).Value = 'Normal'
$ws.Cells.Item(52, 4).Value = 5.22
$ws.Cells.Item(53, 1).Value = 'Jarad, Golgari Lich Lord'
$ws.Cells.Item(53, 2).Value = 'Duel Decks: Izzet vs Golgari'
$ws.Cells.Item(53, 3).Value = 'Foil'
$ws.Cells.Item(53, 4).Value = 1.69
$ws.Cells.Item(54, 1).Value = 'The Gitrog Monster'
$ws.Cells.Item(54, 2).Value = 'Shadows over Innistrad'
$ws.Cells.Item(54, 3).Value = 'Normal'
$ws.Cells.Item(54, 4).Value = 3.75
$ws.Cells.Item(55, 1).Value = 'Arlinn Kord: Arlinn, Embraced by the Moon'
$ws.Cells.Item(55, 2).Value = 'Shadows over Innistrad'
$ws.Cells.Item(55, 3).Value = 'Normal'
$ws.Cells.Item(55, 4).Value = 4.32
$ws.Cells.Item(56, 1).Value = 'Rith, Liberated Primeval'
$ws.Cells.Item(56, 2).Value = 'Dominaria United: Extras'
$ws.Cells.Item(56, 3).Value = 'V.2'
$ws.Cells.Item(56, 4).Value = 2.79
$ws.Cells.Item(57, 1).Value = 'Miirym, Sentinel Wyrm'
$ws.Cells.Item(57, 2).Value = 'Commander Legends: Battle For Baldur''s Gate'
$ws.Cells.Item(57, 3).Value = 'Normal'
$ws.Cells.Item(57, 4).Value = 1.48
$ws.Cells.Item(58, 1).Value = 'Aether Vial'
$ws.Cells.Item(58, 2).Value = 'Double Masters 2022'
$ws.Cells.Item(58, 3).Value = 'Normal'
$ws.Cells.Item(58, 4).Value = 11.89
$ws.Cells.Item(59, 1).Value = 'Akroma''s Memorial'
$ws.Cells.Item(59, 2).Value = 'Magic 2013'
$ws.Cells.Item(59, 3).Value = 'Normal'
$ws.Cells.Item(59, 4).Value = 7.22
$ws.Cells.Item(60, 1).Value = 'Amulet of Vigor'
$ws.Cells.Item(60, 2).Value = 'Mystery Booster'
$ws.Cells.Item(60, 3).Value = 'Foil'
$ws.Cells.Item(60, 4).Value = 17.82
$ws.Cells.Item(61, 1).Value = 'Astral Cornucopia'
$ws.Cells.Item(61, 2).Value = 'Born of the Gods'
$ws.Cells.Item(61, 3).Value = 'Normal'
$ws.Cells.Item(61, 4).Value = 1.15
$ws.Cells.Item(62, 1).Value = 'Coat of Arms'
$ws.Cells.Item(62, 2).Value = 'Mystery Booster'
$ws.Cells.Item(62, 3).Value = 'Normal'
$ws.Cells.Item(62, 4).Value = 10.17
$ws.Cells.Item(63, 1).Value = 'Dolmen Gate'
$ws.Cells.Item(63, 2).Value = 'Mystery Booster'
$ws.Cells.Item(63, 3).Value = 'Normal'
$ws.Cells.Item(63, 4).Value = 9
$ws.Cells.Item(64, 1).Value = 'Elbrus, the Binding Blade: Withengar Unbound'
$ws.Cells.Item(64, 2).Value = 'Dark Ascension'
$ws.Cells.Item(64, 3).Value = 'Normal'
$ws.Cells.Item(64, 4).Value = 3.64
$ws.Cells.Item(65, 1).Value = 'Grafdigger''s Cage'
$ws.Cells.Item(65, 2).Value = 'Dark Ascension'
$ws.Cells.Item(65, 3).Value = 'Normal'
$ws.Cells.Item(65, 4).Value = 1.97
$ws.Cells.Item(66, 1).Value = 'Helm of Awakening'
$ws.Cells.Item(66, 2).Value = 'Dominaria Remastered'
$ws.Cells.Item(66, 3).Value = 'Foil'
$ws.Cells.Item(66, 4).Value = 1.75
$ws.Cells.Item(67, 1).Value = 'Illusionist''s Bracers'
$ws.Cells.Item(67, 2).Value = 'Gatecrash'
$ws.Cells.Item(67, 3).Value = 'Normal'
$ws.Cells.Item(67, 4).Value = 6.05
$ws.Cells.Item(68, 1).Value = 'Karn, Living Legacy'
$ws.Cells.Item(68, 2).Value = 'Dominaria United Promos'
$ws.Cells.Item(68, 3).Value = 'V.2'
$ws.Cells.Item(68, 4).Value = 3.8
$ws.Cells.Item(69, 1).Value = 'Swiftfoot Boots'
$ws.Cells.Item(69, 2).Value = 'Commander 2017'
$ws.Cells.Item(69, 3).Value = 'Normal'
$ws.Cells.Item(69, 4).Value = 1.01

# Restore selection to match the authored state
$ws.Range("E69").Select()

Write-Host "done"
